$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty row 14 (date, start/end time, and
# source/project/task/description) - like a new journal entry.
$ws.Range("B14").Value = 43907
$ws.Range("C14").Value = 0.58333333333333337
$ws.Range("D14").Value = 0.61319444444444449
$ws.Range("F14").Value = "CLion"
$ws.Range("G14").Value = "Bataille Navale"
$ws.Range("H14").Value = "Programmation du jeu"
$ws.Range("I14").Value = "Création des différentes fonctions"

# Move the active selection from D14 to H14
$ws.Range("H14").Select()
